$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Properties Table" sheet: update the CREATE TABLE documentation string
#    in H1 (Flags INTEGER -> Flags VARCHAR(3)). Changing this shared string
#    automatically drops the now-unused old copy and appends the new one at
#    the end of the shared-string table, which is what re-numbers every
#    other <v> reference on the "Field Explanation" / "Types" sheets.
# ---------------------------------------------------------------------------
$wsProps = $wb.Worksheets.Item("Properties Table")

$wsProps.Range("H1").Value = "CREATE TABLE [Properties] ([ID] INTEGER  NOT NULL PRIMARY KEY AUTOINCREMENT,[Name] VARCHAR(30)  UNIQUE NOT NULL,[Location] VARCHAR(150)  UNIQUE NOT NULL,[Flags] VARCHAR(3)  NOT NULL,[Staff] INTEGER DEFAULT '0' NOT NULL,[StaffCap] INTEGER DEFAULT '10' NOT NULL,[Cost] INTEGER  NOT NULL, [TypeID] INTEGER  NOT NULL)"

# 2) H2's formula now wraps the Flags value (column C) in quotes so it is
#    inserted as text instead of a bare integer.
$wsProps.Range("H2").Formula = '="INSERT OR IGNORE INTO Properties (Name, Location, Flags, Staff, StaffCap, Cost, Type) VALUES (''" &A2&"'',''" &B2& "'',''"&C2&"'',"&D2&","&E2&","&F2&","&G2&");"'

# 3) Drop the per-row SQL helper formulas in H3:H11 entirely (the shared
#    formula group anchored at H3 is removed, not just cleared of content).
$wsProps.Range("H3:H11").Clear()

# 4) View-state: the active window now has column H scrolled into view and
#    the selection sits on H2 (inside the frozen bottom-left pane).
$wsProps.Activate() | Out-Null
$wsProps.Range("H2").Select() | Out-Null

# ---------------------------------------------------------------------------
# "Field Explanation" sheet: only the on-screen selection moved (E9 -> C9).
# ---------------------------------------------------------------------------
$wsField = $wb.Worksheets.Item("Field Explanation")
$wsField.Activate() | Out-Null
$wsField.Range("C9").Select() | Out-Null

# Leave the workbook focused back on "Properties Table" (the originally
# active/tab-selected sheet).
$wsProps.Activate() | Out-Null
